# Generate Report for Archive
# The localization status moved on from "Ready for handoff" to "In Translation"
# for every tracked file, across the Overview sheet (zh-cn / de-de status
# columns) and each per-locale report sheet (Status column). Updating the
# text makes the status columns narrower, so their widths are refreshed too.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 13.4101845877511

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Range("E1:F1").ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $newWidth
